$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D26").Value = "생성 모델의 새로운 흐름 확산 모델(Diffusion model)에 관하여"

$ws.Range("D37").Value = "[Paper Review] CFLOW-AD: Real-Time Unsupervised Anomaly Detection with Localization via Conditional Normalizing Flows"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1977&mod=document&pageid=1"

$ws.Range("D51").Value = "구글링할 때 2022년 이후의 자료만 검색되게 하고 싶다면?"
$ws.Range("E51").Value = "https://bskyvision.com/1268"

$ws.Range("D52").Value = "[파이썬] R과 차이점 01"
$ws.Range("E52").Value = "http://ds.sumeun.org/?p=2596&utm_source=rss&utm_medium=rss&utm_campaign=%ed%8c%8c%ec%9d%b4%ec%8d%ac-r%ea%b3%bc-%ec%b0%a8%ec%9d%b4%ec%a0%90-01"
